# Remove the "LITERATURE REVIEW" slide (SlideID 277) from the presentation.
$p = $ppt.ActivePresentation

for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $s = $p.Slides.Item($i)
    if ($s.SlideID -eq 277) {
        $s.Delete()
        break
    }
}
